$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.940.63'
$ws.Range('E2').Value = '  +3.28%  '

# Row 3
$ws.Range('D3').Value = '2.974.04'
$ws.Range('E3').Value = '  +0.92%  '

# Row 4
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').Value = "'575.30"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.38%  '

# Row 6
$ws.Range('D6').Value = "'160.70"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +6.92%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').Value = "'0.515"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.80%  '

# Row 9
$ws.Range('D9').Value = '2.969.46'
$ws.Range('E9').Value = '  +0.92%  '

# Row 10
$ws.Range('D10').Value = "'6.70"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.92%  '

# Row 11
$ws.Range('E11').Value = '  +0.52%  '

# Row 12
$ws.Range('D12').Value = "'0.452"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.41%  '

# Row 13
$ws.Range('D13').Value = "'0.0000246"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.01%  '

# Row 14
$ws.Range('D14').Value = "'34.24"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.76%  '

# Row 15
$ws.Range('E15').Value = '  -0.62%  '

# Row 16
$ws.Range('D16').Value = '65.945.81'
$ws.Range('E16').Value = '  +3.52%  '

# Row 17
$ws.Range('D17').Value = '3.469.67'
$ws.Range('E17').Value = '  +1.06%  '

# Row 18
$ws.Range('D18').Value = "'6.86"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.46%  '

# Row 19
$ws.Range('D19').Value = '2.978.33'
$ws.Range('E19').Value = '  +0.94%  '

# Row 20
$ws.Range('D20').Value = "'450.42"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.27%  '

# Row 21
$ws.Range('D21').Value = "'13.77"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.86%  '

# Row 22
$ws.Range('D22').Value = "'0.678"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.24%  '

# Row 23
$ws.Range('D23').Value = "'7.26"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.01%  '

# Row 24
$ws.Range('D24').Value = "'81.82"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.36%  '

# Row 25
$ws.Range('D25').Value = "'2.23"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.67%  '

# Row 26
$ws.Range('D26').Value = "'12.17"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.43%  '

# Row 27
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = "'10.01"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -7.03%  '

# Row 28
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = "'1.00"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.08%  '

# Row 29
$ws.Range('D29').Value = "'8.11"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.37%  '

# Row 30
$ws.Range('E30').Value = '  +12.36%  '

# Row 31
$ws.Range('D31').Value = "'2.59"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.98%  '

# Row 32
$ws.Range('E32').Value = '  -4.74%  '

# Row 33
$ws.Range('D33').Value = "'27.08"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.05%  '

# Row 34
$ws.Range('E34').Value = '  +1.44%  '

# Row 35
$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.09%  '

# Row 36
$ws.Range('D36').Value = "'0.980"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.39%  '

# Row 37
$ws.Range('E37').Value = '  +3.66%  '

# Row 38
$ws.Range('E38').Value = '  -3.67%  '

# Row 39
$ws.Range('D39').Value = "'49.47"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.87%  '

# Row 40
$ws.Range('D40').Value = "'43.56"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.60%  '

# Row 41
$ws.Range('D41').Value = "'2.85"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.08%  '

# Row 42
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = "'0.298"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.77%  '

# Row 43
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = "'0.119"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.54%  '

# Row 44
$ws.Range('D44').Value = "'8.37"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.75%  '

# Row 45
$ws.Range('D45').Value = "'387.09"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.92%  '

# Row 46
$ws.Range('D46').Value = "'0.0355"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.77%  '

# Row 47
$ws.Range('D47').Value = '2.727.43'
$ws.Range('E47').Value = '  -0.36%  '

# Row 48
$ws.Range('D48').Value = "'132.47"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.10%  '

# Row 49
$ws.Range('E49').Value = '  +0.04%  '

# Row 50
$ws.Range('E50').Value = '  +1.17%  '

# Row 51
$ws.Range('D51').Value = "'23.14"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.86%  '
